# Update TPM-derived values in the LR-pairs sheet (Hc-C5ar2.xlsx)
# The original sheet had 6 data rows (rows 2-7); the updated data only has
# 4 data rows (rows 2-5). Row 4 (FAPs->MuSCs) and row 6 (MuSCs->FAPs) of the
# old layout are removed, row ordering is adjusted, and many numeric values
# are refreshed with newly computed TPM-based statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that are dropped from the updated dataset
# (old row 7 = MuSCs/Hc/C5ar2/MuSCs, old row 6 = MuSCs/Hc/C5ar2/FAPs).
# Delete from the bottom up so row indices of earlier rows stay valid.
$ws.Rows(7).Delete()
$ws.Rows(6).Delete()

# ---- Row 2: FAPs -> ECs (unchanged pairing, refreshed stats) ----
$ws.Range("I2").Value = 0.8736649195182647
$ws.Range("J2").Value = 0.8736649195182647
$ws.Range("M2").Value = 0.008059999999999999
$ws.Range("N2").Value = 0.02418
$ws.Range("O2").Value = 0.1431624817198444
$ws.Range("P2").Value = 0.1431624817198444
$ws.Range("Q2").Value = 0.002397651186666666
$ws.Range("R2").Value = 0.02157886068
$ws.Range("S2").Value = 0.1250760380698029
$ws.Range("T2").Value = 0.1250760380698029

# ---- Row 3: FAPs -> FAPs (unchanged pairing, refreshed stats) ----
$ws.Range("I3").Value = 0.8736649195182647
$ws.Range("J3").Value = 0.8736649195182647
$ws.Range("O3").Value = 0.8568375182801556
$ws.Range("P3").Value = 0.8568375182801555
$ws.Range("S3").Value = 0.7485888814484619
$ws.Range("T3").Value = 0.7485888814484618

# ---- Row 4: now MuSCs -> ECs (was FAPs -> MuSCs) ----
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "ECs"
$ws.Range("G4").Value = 0.043016
$ws.Range("H4").Value = 0.129048
$ws.Range("I4").Value = 0.1263350804817352
$ws.Range("J4").Value = 0.1263350804817352
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.008059999999999999
$ws.Range("N4").Value = 0.02418
$ws.Range("O4").Value = 0.1431624817198444
$ws.Range("P4").Value = 0.1431624817198444
$ws.Range("Q4").Value = 0.0003467089599999999
$ws.Range("R4").Value = 0.00312038064
$ws.Range("S4").Value = 0.01808644365004149
$ws.Range("T4").Value = 0.01808644365004149

# ---- Row 5: now MuSCs -> FAPs (was MuSCs -> ECs) ----
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.043016
$ws.Range("H5").Value = 0.129048
$ws.Range("I5").Value = 0.1263350804817352
$ws.Range("J5").Value = 0.1263350804817352
$ws.Range("M5").Value = 0.04823966666666667
$ws.Range("N5").Value = 0.144719
$ws.Range("O5").Value = 0.8568375182801556
$ws.Range("P5").Value = 0.8568375182801555
$ws.Range("Q5").Value = 0.002075077501333334
$ws.Range("R5").Value = 0.018675697512
$ws.Range("S5").Value = 0.1082486368316937
$ws.Range("T5").Value = 0.1082486368316937
